# Update column F ("想去人数") values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Row -> new value for column F, shared by both sheets except row 4
# (展览 goes 2020->2157, 全部类型 goes 2021->2157; both land on 2157)
$updates = @{
    2  = 838
    4  = 2157
    5  = 53
    6  = 12462
    9  = 500
    10 = 449
    11 = 1140
    12 = 927
    13 = 13628
    14 = 13882
    16 = 164
    18 = 41
    19 = 1039
    20 = 105
    22 = 368
    23 = 4983
    24 = 234
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
